$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PUBLIC.API")

# Add the new "Unhandled exceptions" error-code row right after the existing data (row 29 -> new row 31)
$ws.Range("A31").Value = "API"
$ws.Range("B31").Value = "Global"
$ws.Range("C31").Value = "PUBLIC_API_99999"
$ws.Range("D31").Value = "Unhandled exceptions"
$ws.Range("E31").Value = "null"

# Resize the existing table (ListObject) so it covers the new row
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("A1:E31"))

# Update selection to mirror the authored state after entering the new data
$ws.Range("C31").Select()
